$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 22-23 need the Fecha column to use the same date number format as the rest of column D
$ws.Range("D22:D23").NumberFormat = $ws.Range("D2").NumberFormat

# Row 2
$ws.Range("D2").Value = 44553
$ws.Range("K2").Value = 'Modesto'
$ws.Range("N2").Value = 23000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 23500
$ws.Range("Q2").Value = '$/caja 16 kilos'
$ws.Range("S2").Value = 1469
$ws.Range("T2").Value = 16
# Row 3
$ws.Range("D3").Value = 44553
$ws.Range("K3").Value = 'Modesto'
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21500
$ws.Range("Q3").Value = '$/caja 16 kilos'
$ws.Range("S3").Value = 1344
$ws.Range("T3").Value = 16
# Row 4
$ws.Range("D4").Value = 44553
$ws.Range("K4").Value = 'Modesto'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 240
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = '$/caja 16 kilos'
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 1094
$ws.Range("T4").Value = 16
# Row 5
$ws.Range("D5").Value = 44566
$ws.Range("K5").Value = 'Modesto'
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 23000
$ws.Range("O5").Value = 24000
$ws.Range("P5").Value = 23500
$ws.Range("S5").Value = 1306
# Row 6
$ws.Range("D6").Value = 44566
$ws.Range("K6").Value = 'Modesto'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 1194
# Row 7
$ws.Range("D7").Value = 44559
$ws.Range("K7").Value = 'Modesto'
$ws.Range("M7").Value = 400
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 26000
$ws.Range("P7").Value = 25500
$ws.Range("S7").Value = 1417
# Row 8
$ws.Range("D8").Value = 44559
$ws.Range("K8").Value = 'Modesto'
$ws.Range("M8").Value = 320
$ws.Range("N8").Value = 22000
$ws.Range("O8").Value = 23000
$ws.Range("P8").Value = 22500
$ws.Range("S8").Value = 1250
# Row 9
$ws.Range("D9").Value = 44189
$ws.Range("K9").Value = 'Dina'
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 23500
$ws.Range("O9").Value = 24000
$ws.Range("P9").Value = 23750
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1319
# Row 10
$ws.Range("D10").Value = 44189
$ws.Range("K10").Value = 'Dina'
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 21500
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 21750
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1208
# Row 11
$ws.Range("D11").Value = 44175
$ws.Range("K11").Value = 'Castle Brite'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 21000
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 21500
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 1194
# Row 12
$ws.Range("D12").Value = 44546
$ws.Range("K12").Value = 'Castle Brite'
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 22500
$ws.Range("O12").Value = 23000
$ws.Range("P12").Value = 22750
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 1264
# Row 13
$ws.Range("D13").Value = 44546
$ws.Range("K13").Value = 'Castle Brite'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 20500
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20750
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("S13").Value = 1153
$ws.Range("T13").Value = 18
# Row 14
$ws.Range("D14").Value = 44573
$ws.Range("L14").Value = 'Especial'
$ws.Range("N14").Value = 20500
$ws.Range("O14").Value = 21000
$ws.Range("P14").Value = 20750
$ws.Range("Q14").Value = '$/caja 18 kilos'
$ws.Range("S14").Value = 1153
$ws.Range("T14").Value = 18
# Row 15
$ws.Range("D15").Value = 44573
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 400
$ws.Range("N15").Value = 17500
$ws.Range("P15").Value = 17750
$ws.Range("Q15").Value = '$/caja 18 kilos'
$ws.Range("S15").Value = 986
$ws.Range("T15").Value = 18
# Row 19
$ws.Range("D19").Value = 44552
$ws.Range("K19").Value = 'Castle Brite'
$ws.Range("M19").Value = 360
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 21000
$ws.Range("P19").Value = 20500
$ws.Range("R19").Value = 'Región Metropolitana'
$ws.Range("S19").Value = 1139
# Row 20
$ws.Range("D20").Value = 44552
$ws.Range("K20").Value = 'Castle Brite'
$ws.Range("M20").Value = 280
$ws.Range("N20").Value = 18000
$ws.Range("O20").Value = 19000
$ws.Range("P20").Value = 18500
$ws.Range("R20").Value = 'Región Metropolitana'
$ws.Range("S20").Value = 1028
# Row 21
$ws.Range("D21").Value = 44545
$ws.Range("L21").Value = 'Especial'
$ws.Range("M21").Value = 340
$ws.Range("N21").Value = 22500
$ws.Range("O21").Value = 23000
$ws.Range("P21").Value = 22750
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 1264
# Row 22
$ws.Range("A22").Value = 2
$ws.Range("B22").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C22").Value = 'Coquimbo'
$ws.Range("D22").Value = 44545
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 'Fruta'
$ws.Range("G22").Value = 100103
$ws.Range("H22").Value = 'Frutos de hueso (carozo)'
$ws.Range("I22").Value = 100103003
$ws.Range("J22").Value = 'Damasco'
$ws.Range("K22").Value = 'Castle Brite'
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 400
$ws.Range("N22").Value = 20500
$ws.Range("O22").Value = 21000
$ws.Range("P22").Value = 20750
$ws.Range("Q22").Value = '$/caja 18 kilos'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1153
$ws.Range("T22").Value = 18
# Row 23
$ws.Range("A23").Value = 2
$ws.Range("B23").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C23").Value = 'Coquimbo'
$ws.Range("D23").Value = 44545
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = 'Fruta'
$ws.Range("G23").Value = 100103
$ws.Range("H23").Value = 'Frutos de hueso (carozo)'
$ws.Range("I23").Value = 100103003
$ws.Range("J23").Value = 'Damasco'
$ws.Range("K23").Value = 'Castle Brite'
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 300
$ws.Range("N23").Value = 15500
$ws.Range("O23").Value = 16000
$ws.Range("P23").Value = 15750
$ws.Range("Q23").Value = '$/caja 18 kilos'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 875
$ws.Range("T23").Value = 18
